$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply integer number format (numFmtId 1 => "0") to the data columns C:F, rows 2-62
$ws.Range("C2:F62").NumberFormat = "0"

# Widen column B to fit the study names (target stored width 25.42578125;
# the engine quantizes ColumnWidth to steps of 1/6, so 24.6667 is the input
# that lands on the closest achievable stored width, 25.5)
$ws.Columns.Item(2).ColumnWidth = 24.666666666666668

# Update the active selection cell to H24
$ws.Range("H24").Select()
